$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

$rows = @(
    @(4, "TPE", "J", "P", 10000),
    @(5, "TPE", "E", "H1", 4000),
    @(6, "TPE", "V", "H2", 2400)
)

$r = 5
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

$ws.Range("E7").Select() | Out-Null
